$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, [string]$val) {
    # Force text interpretation so numeric-looking strings (e.g. "1.000",
    # "0.9998") are stored as literal text rather than being parsed into
    # numbers/dates, then restore the cell's original style so no visible
    # formatting change is introduced.
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "25.084.50"
Set-TextValue $ws.Range("E2") "  -3.01%  "
Set-TextValue $ws.Range("D3") "1.649.86"
Set-TextValue $ws.Range("E3") "  -4.97%  "
Set-TextValue $ws.Range("D4") "0.9998"
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "237.14"
Set-TextValue $ws.Range("E5") "  -2.16%  "
Set-TextValue $ws.Range("E6") "  +0.01%  "
Set-TextValue $ws.Range("D7") "0.4787"
Set-TextValue $ws.Range("E7") "  -7.78%  "
Set-TextValue $ws.Range("D8") "0.2626"
Set-TextValue $ws.Range("E8") "  -4.07%  "
Set-TextValue $ws.Range("D9") "0.06039"
Set-TextValue $ws.Range("E9") "  -1.98%  "
Set-TextValue $ws.Range("D10") "0.07097"
Set-TextValue $ws.Range("E10") "  -1.27%  "
Set-TextValue $ws.Range("D11") "1.647.73"
Set-TextValue $ws.Range("E11") "  -5.14%  "
Set-TextValue $ws.Range("D12") "14.47"
Set-TextValue $ws.Range("E12") "  -3.25%  "
Set-TextValue $ws.Range("D13") "0.6187"
Set-TextValue $ws.Range("E13") "  -3.73%  "
Set-TextValue $ws.Range("D14") "4.568"
Set-TextValue $ws.Range("E14") "  -0.96%  "
Set-TextValue $ws.Range("D15") "73.10"
Set-TextValue $ws.Range("E15") "  -5.37%  "
Set-TextValue $ws.Range("D16") "1.000"
Set-TextValue $ws.Range("E16") "  +0.04%  "
Set-TextValue $ws.Range("D17") "0.9994"
Set-TextValue $ws.Range("E17") "  -0.06%  "
Set-TextValue $ws.Range("D18") "25.068.01"
Set-TextValue $ws.Range("E18") "  -3.21%  "
Set-TextValue $ws.Range("D19") "11.37"
Set-TextValue $ws.Range("E19") "  -3.34%  "
Set-TextValue $ws.Range("D20") "0.000006566"
Set-TextValue $ws.Range("E20") "  -3.20%  "
Set-TextValue $ws.Range("D21") "4.416"
Set-TextValue $ws.Range("E21") "  +3.15%  "
Set-TextValue $ws.Range("D22") "1.864.14"
Set-TextValue $ws.Range("E22") "  -5.01%  "
Set-TextValue $ws.Range("D23") "8.473"
Set-TextValue $ws.Range("E23") "  -1.86%  "
Set-TextValue $ws.Range("D24") "5.239"
Set-TextValue $ws.Range("E24") "  -0.98%  "
Set-TextValue $ws.Range("D25") "133.88"
Set-TextValue $ws.Range("E25") "  -2.26%  "
Set-TextValue $ws.Range("D26") "14.73"
Set-TextValue $ws.Range("E26") "  -3.18%  "
Set-TextValue $ws.Range("D27") "1.394"
Set-TextValue $ws.Range("E27") "  -7.72%  "
Set-TextValue $ws.Range("D28") "1.693"
Set-TextValue $ws.Range("E28") "  -4.45%  "
Set-TextValue $ws.Range("D29") "101.93"
Set-TextValue $ws.Range("E29") "  -3.12%  "
Set-TextValue $ws.Range("D30") "3.791"
Set-TextValue $ws.Range("E30") "  -4.43%  "
Set-TextValue $ws.Range("E31") "  -3.91%  "
Set-TextValue $ws.Range("E32") "  -2.59%  "
Set-TextValue $ws.Range("D33") "0.04544"
Set-TextValue $ws.Range("E33") "  -2.78%  "
Set-TextValue $ws.Range("D34") "2.608"
Set-TextValue $ws.Range("E34") "  -1.65%  "
Set-TextValue $ws.Range("D35") "0.9409"
Set-TextValue $ws.Range("E35") "  -4.96%  "
Set-TextValue $ws.Range("D36") "0.5791"
Set-TextValue $ws.Range("E36") "  -6.47%  "
Set-TextValue $ws.Range("D37") "2.627"
Set-TextValue $ws.Range("E37") "  -2.29%  "
Set-TextValue $ws.Range("D38") "0.01538"
Set-TextValue $ws.Range("E39") "  +12.41%  "
Set-TextValue $ws.Range("D40") "1.000"
Set-TextValue $ws.Range("E40") "  +0.05%  "
Set-TextValue $ws.Range("E41") "  -5.38%  "
Set-TextValue $ws.Range("D42") "98.73"
Set-TextValue $ws.Range("E42") "  -0.98%  "
Set-TextValue $ws.Range("D43") "0.3701"
Set-TextValue $ws.Range("E43") "  -3.98%  "
Set-TextValue $ws.Range("D44") "4.790"
Set-TextValue $ws.Range("E44") "  -4.35%  "
Set-TextValue $ws.Range("D45") "0.1124"
Set-TextValue $ws.Range("E45") "  -0.35%  "
Set-TextValue $ws.Range("D46") "6.031"
Set-TextValue $ws.Range("E46") "  -3.48%  "
Set-TextValue $ws.Range("D47") "0.05157"
Set-TextValue $ws.Range("E47") "  -1.09%  "
Set-TextValue $ws.Range("D48") "52.12"
Set-TextValue $ws.Range("E48") "  -5.28%  "
Set-TextValue $ws.Range("D49") "29.50"
Set-TextValue $ws.Range("E49") "  -3.60%  "
Set-TextValue $ws.Range("D50") "0.9998"
Set-TextValue $ws.Range("E50") "  -0.07%  "
Set-TextValue $ws.Range("D51") "0.3333"
Set-TextValue $ws.Range("E51") "  -2.39%  "
